$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (tab name) to reflect the updated "through" date
$ws.Name = "Through 2022-09-09"

# Update the label cell for the September row
$ws.Range("A10").Value = "September (through 09-09)"

# Update September row (row 10) values for columns C..I (B10 unchanged)
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = 22
$ws.Range("E10").Value = 14
$ws.Range("F10").Value = 21
$ws.Range("G10").Value = 28
$ws.Range("H10").Value = 38
$ws.Range("I10").Value = 39

# Update Total row (row 11) values for columns C..I (B11 unchanged)
$ws.Range("C11").Value = 396
$ws.Range("D11").Value = 573
$ws.Range("E11").Value = 504
$ws.Range("F11").Value = 370
$ws.Range("G11").Value = 812
$ws.Range("H11").Value = 1108
$ws.Range("I11").Value = 1176
